# Auto-update timeframe-hour MOB report (reporthouronfcst)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 852
$ws.Range("D2").Value = 867
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 14
$ws.Range("I2").Value = 962
$ws.Range("J2").Value = -9.875259875259879

$ws.Range("C3").Value = 377
$ws.Range("D3").Value = 377
$ws.Range("G3").Value = 11
$ws.Range("I3").Value = 520
$ws.Range("J3").Value = -27.5

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = -50

$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 244
$ws.Range("D5").Value = 306
$ws.Range("E5").Value = 46
$ws.Range("F5").Value = 11
$ws.Range("I5").Value = 252
$ws.Range("J5").Value = 21.42857142857142

$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 56
$ws.Range("I6").Value = 60
$ws.Range("J6").Value = -6.666666666666665

$ws.Range("C7").Value = 58
$ws.Range("D7").Value = 61
$ws.Range("E7").Value = 3
$ws.Range("I7").Value = 101
$ws.Range("J7").Value = -39.6039603960396

$ws.Range("C8").Value = 213
$ws.Range("D8").Value = 214
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 6
$ws.Range("J8").Value = 919.047619047619

$ws.Range("C9").Value = 38
$ws.Range("D9").Value = 37
$ws.Range("G9").Value = 2
$ws.Range("J9").Value = -22.91666666666666

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 246
$ws.Range("D10").Value = 325
$ws.Range("E10").Value = 84
$ws.Range("G10").Value = 5
$ws.Range("I10").Value = 537
$ws.Range("J10").Value = -39.47858472998138

$ws.Range("C11").Value = 269
$ws.Range("D11").Value = 272
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 5
$ws.Range("I11").Value = 303
$ws.Range("J11").Value = -10.23102310231023

$ws.Range("C12").Value = 424
$ws.Range("D12").Value = 615
$ws.Range("E12").Value = 38
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 147
$ws.Range("I12").Value = 760.5
$ws.Range("J12").Value = -19.13214990138067

$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 27
$ws.Range("I13").Value = 733
$ws.Range("J13").Value = -96.31650750341065

$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 529
$ws.Range("D14").Value = 627
$ws.Range("E14").Value = 25
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = 73
$ws.Range("I14").Value = 860
$ws.Range("J14").Value = -27.09302325581395

$ws.Range("C15").Value = 195
$ws.Range("D15").Value = 203
$ws.Range("E15").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("I15").Value = 278
$ws.Range("J15").Value = -26.97841726618705

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 63
$ws.Range("D17").Value = 70
$ws.Range("E17").Value = 6
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = -42.62295081967213

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("J18").Value = -50

$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 1
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 14.28571428571428

$ws.Range("C20").Value = 43
$ws.Range("D20").Value = 48
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 1
$ws.Range("I20").Value = 69
$ws.Range("J20").Value = -30.43478260869566

